# Generate Report for Handback
# - Updates the Overview status text to reflect that the handback is
#   complete and in sync with en-US.
# - Populates the per-language (zh-cn / de-de) sheets with the latest
#   target/handback file names and the handback timestamp, and links the
#   "Latest Target File" cell back to the source markdown file.
# - Widens a few columns so the new, longer text fits.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

$mdFileName  = "17ae3688-f601-4fda-9ec4-75a4b1907617.md"
$mdUrl       = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d392c61bfe86506dae9e8e5500e5f45b7e36025d/e2e/17ae3688-f601-4fda-9ec4-75a4b1907617.md"

$statusText  = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: status columns (zh-cn / de-de) now read "Handed back"
# ---------------------------------------------------------------------
$ws_overview.Range("E2").Value = $statusText
$ws_overview.Range("F2").Value = $statusText

$ws_overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$ws_overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws_zhcn.Range("C2").Value = $statusText

$ws_zhcn.Range("I2").Value = $mdFileName
$ws_zhcn.Hyperlinks.Add($ws_zhcn.Range("I2"), $mdUrl, "", "", $mdFileName) | Out-Null

$ws_zhcn.Range("J2").Value = "17ae3688-f601-4fda-9ec4-75a4b1907617.e388192d03ef0a6dd716f17b66d359b14047db4d.zh-cn.xlf"
$ws_zhcn.Range("K2").Value = "2016-09-01 23:04:16"

$ws_zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$ws_zhcn.Columns.Item(9).ColumnWidth = 40
$ws_zhcn.Columns.Item(10).ColumnWidth = 40

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws_dede.Range("C2").Value = $statusText

$ws_dede.Range("I2").Value = $mdFileName
$ws_dede.Hyperlinks.Add($ws_dede.Range("I2"), $mdUrl, "", "", $mdFileName) | Out-Null

$ws_dede.Range("J2").Value = "17ae3688-f601-4fda-9ec4-75a4b1907617.e388192d03ef0a6dd716f17b66d359b14047db4d.de-de.xlf"
$ws_dede.Range("K2").Value = "2016-09-01 23:04:23"

$ws_dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$ws_dede.Columns.Item(9).ColumnWidth = 40
$ws_dede.Columns.Item(10).ColumnWidth = 40
